$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 15).Value = 1.11
$ws.Cells.Item(2, 16).Value = 6.5
$ws.Cells.Item(2, 17).Value = 1.4
$ws.Cells.Item(2, 18).Value = 3
$ws.Cells.Item(2, 23).Value = 1.22
$ws.Cells.Item(2, 24).Value = 4
$ws.Cells.Item(2, 25).Value = 1.67
$ws.Cells.Item(2, 26).Value = 2.1
$ws.Cells.Item(2, 27).Value = 12
$ws.Cells.Item(2, 28).Value = 9
$ws.Cells.Item(2, 29).Value = 9.5
$ws.Cells.Item(2, 30).Value = 10
$ws.Cells.Item(2, 32).Value = 23
$ws.Cells.Item(2, 33).Value = 21
$ws.Cells.Item(2, 34).Value = 12
$ws.Cells.Item(2, 36).Value = 51
$ws.Cells.Item(2, 37).Value = 201
$ws.Cells.Item(2, 42).Value = 67
$ws.Cells.Item(3, 7).Value = 4.1
$ws.Cells.Item(3, 17).Value = 1.91
$ws.Cells.Item(3, 18).Value = 1.99
$ws.Cells.Item(4, 7).Value = 1.48
$ws.Cells.Item(4, 9).Value = 6.25
$ws.Cells.Item(4, 30).Value = 10
$ws.Cells.Item(4, 39).Value = 34
$ws.Cells.Item(4, 40).Value = 19
$ws.Cells.Item(5, 17).Value = 2.5
$ws.Cells.Item(5, 18).Value = 1.5
$ws.Cells.Item(5, 19).Value = 4
$ws.Cells.Item(5, 20).Value = 1.25
$ws.Cells.Item(5, 21).Value = 5
$ws.Cells.Item(5, 22).Value = 1.17
$ws.Cells.Item(5, 44).Value = 1.88
$ws.Cells.Item(5, 45).Value = 1.98
$ws.Cells.Item(6, 7).Value = 2.3
$ws.Cells.Item(6, 9).Value = 3.2
$ws.Cells.Item(6, 10).Value = 3
$ws.Cells.Item(6, 14).Value = 8.5
$ws.Cells.Item(6, 30).Value = 21
$ws.Cells.Item(6, 31).Value = 19
$ws.Cells.Item(6, 40).Value = 12
$ws.Cells.Item(7, 13).Value = 1.13
$ws.Cells.Item(7, 14).Value = 6
$ws.Cells.Item(7, 19).Value = 5
$ws.Cells.Item(8, 9).Value = 4.33
$ws.Cells.Item(8, 17).Value = 3.5
$ws.Cells.Item(8, 18).Value = 1.3
$ws.Cells.Item(8, 23).Value = 1.78
$ws.Cells.Item(8, 24).Value = 2.03
$ws.Cells.Item(8, 27).Value = 4.5
$ws.Cells.Item(8, 42).Value = 51
$ws.Cells.Item(10, 15).Value = 1.44
$ws.Cells.Item(10, 16).Value = 2.63
$ws.Cells.Item(10, 19).Value = 3.95
$ws.Cells.Item(10, 20).Value = 1.24
$ws.Cells.Item(10, 44).Value = 1.88
$ws.Cells.Item(10, 45).Value = 1.98
$ws.Cells.Item(13, 7).Value = 4.8
$ws.Cells.Item(13, 8).Value = 3.85
$ws.Cells.Item(13, 9).Value = 1.65
$ws.Cells.Item(13, 10).Value = 5
$ws.Cells.Item(13, 12).Value = 2.22
$ws.Cells.Item(13, 13).Value = 1.05
$ws.Cells.Item(13, 14).Value = 8.25
$ws.Cells.Item(13, 16).Value = 3.5
$ws.Cells.Item(13, 17).Value = 1.82
$ws.Cells.Item(13, 18).Value = 1.93
$ws.Cells.Item(13, 21).Value = 2.95
$ws.Cells.Item(13, 22).Value = 1.37
$ws.Cells.Item(13, 25).Value = 1.8
$ws.Cells.Item(13, 26).Value = 1.9
$ws.Cells.Item(13, 28).Value = 30
$ws.Cells.Item(13, 29).Value = 16.5
$ws.Cells.Item(13, 32).Value = 55
$ws.Cells.Item(13, 33).Value = 8.25
$ws.Cells.Item(13, 35).Value = 17
$ws.Cells.Item(13, 38).Value = 6.8
$ws.Cells.Item(13, 40).Value = 8.5
$ws.Cells.Item(13, 42).Value = 14
$ws.Cells.Item(13, 43).Value = 28
$ws.Cells.Item(14, 7).Value = 4.75
$ws.Cells.Item(14, 8).Value = 4
$ws.Cells.Item(14, 9).Value = 1.6
$ws.Cells.Item(14, 10).Value = 5.5
$ws.Cells.Item(14, 11).Value = 2.25
$ws.Cells.Item(14, 12).Value = 2.2
$ws.Cells.Item(14, 15).Value = 1.29
$ws.Cells.Item(14, 16).Value = 3.5
$ws.Cells.Item(14, 17).Value = 1.9
$ws.Cells.Item(14, 18).Value = 1.95
$ws.Cells.Item(14, 21).Value = 3.25
$ws.Cells.Item(14, 22).Value = 1.33
$ws.Cells.Item(14, 23).Value = 1.36
$ws.Cells.Item(14, 24).Value = 3
$ws.Cells.Item(14, 25).Value = 1.91
$ws.Cells.Item(14, 26).Value = 1.91
$ws.Cells.Item(14, 28).Value = 26
$ws.Cells.Item(14, 31).Value = 41
$ws.Cells.Item(14, 33).Value = 12
$ws.Cells.Item(14, 37).Value = 301
$ws.Cells.Item(14, 38).Value = 7
$ws.Cells.Item(14, 39).Value = 7.5
$ws.Cells.Item(14, 41).Value = 12
$ws.Cells.Item(14, 43).Value = 26
$ws.Cells.Item(15, 7).Value = 2.5
$ws.Cells.Item(15, 8).Value = 3.4
$ws.Cells.Item(15, 9).Value = 2.7
$ws.Cells.Item(15, 11).Value = 2.2
$ws.Cells.Item(15, 12).Value = 3.4
$ws.Cells.Item(15, 13).Value = 1.05
$ws.Cells.Item(15, 14).Value = 11
$ws.Cells.Item(15, 15).Value = 1.25
$ws.Cells.Item(15, 16).Value = 4
$ws.Cells.Item(15, 17).Value = 1.85
$ws.Cells.Item(15, 18).Value = 2
$ws.Cells.Item(15, 21).Value = 3.2
$ws.Cells.Item(15, 22).Value = 1.36
$ws.Cells.Item(15, 23).Value = 1.4
$ws.Cells.Item(15, 24).Value = 2.75
$ws.Cells.Item(15, 25).Value = 1.73
$ws.Cells.Item(15, 26).Value = 2
$ws.Cells.Item(15, 27).Value = 9
$ws.Cells.Item(15, 28).Value = 13
$ws.Cells.Item(15, 33).Value = 11
$ws.Cells.Item(15, 35).Value = 13
$ws.Cells.Item(15, 36).Value = 41
$ws.Cells.Item(15, 37).Value = 201
$ws.Cells.Item(15, 38).Value = 9.5
$ws.Cells.Item(15, 40).Value = 10
$ws.Cells.Item(15, 42).Value = 21
$ws.Cells.Item(15, 43).Value = 29
$ws.Cells.Item(16, 7).Value = 1.5
$ws.Cells.Item(16, 8).Value = 4
$ws.Cells.Item(16, 9).Value = 6.5
$ws.Cells.Item(16, 12).Value = 7.5
$ws.Cells.Item(16, 13).Value = 1.07
$ws.Cells.Item(16, 14).Value = 9
$ws.Cells.Item(16, 15).Value = 1.4
$ws.Cells.Item(16, 16).Value = 2.75
$ws.Cells.Item(16, 25).Value = 2.38
$ws.Cells.Item(16, 26).Value = 1.53
$ws.Cells.Item(16, 30).Value = 10
$ws.Cells.Item(16, 40).Value = 21
$ws.Cells.Item(16, 41).Value = 81
$ws.Cells.Item(16, 42).Value = 51
$ws.Cells.Item(20, 7).Value = 1.75
$ws.Cells.Item(20, 9).Value = 4.1
$ws.Cells.Item(20, 10).Value = 2.3
$ws.Cells.Item(20, 12).Value = 4
$ws.Cells.Item(20, 19).Value = 1.8
$ws.Cells.Item(20, 20).Value = 2.05
$ws.Cells.Item(20, 31).Value = 13
$ws.Cells.Item(21, 17).Value = 1.73
$ws.Cells.Item(21, 18).Value = 2.08
$ws.Cells.Item(22, 7).Value = 1.9
$ws.Cells.Item(22, 9).Value = 3.5
$ws.Cells.Item(22, 12).Value = 3.75
$ws.Cells.Item(22, 17).Value = 1.48
$ws.Cells.Item(22, 18).Value = 2.6
$ws.Cells.Item(22, 21).Value = 2.1
$ws.Cells.Item(22, 22).Value = 1.67
$ws.Cells.Item(22, 25).Value = 1.44
$ws.Cells.Item(22, 26).Value = 2.63
$ws.Cells.Item(22, 27).Value = 12
$ws.Cells.Item(22, 30).Value = 19
$ws.Cells.Item(22, 33).Value = 21
$ws.Cells.Item(22, 39).Value = 21
$ws.Cells.Item(22, 42).Value = 23
$ws.Cells.Item(22, 43).Value = 23
$ws.Cells.Item(23, 17).Value = 1.5
$ws.Cells.Item(23, 18).Value = 2.5
$ws.Cells.Item(24, 7).Value = 2.8
$ws.Cells.Item(24, 9).Value = 2.7
$ws.Cells.Item(24, 10).Value = 3.75
$ws.Cells.Item(24, 12).Value = 3.6
$ws.Cells.Item(24, 13).Value = 1.11
$ws.Cells.Item(24, 14).Value = 6.5
$ws.Cells.Item(24, 15).Value = 1.5
$ws.Cells.Item(24, 16).Value = 2.5
$ws.Cells.Item(24, 17).Value = 2.6
$ws.Cells.Item(24, 18).Value = 1.48
$ws.Cells.Item(24, 21).Value = 5.5
$ws.Cells.Item(24, 22).Value = 1.14
$ws.Cells.Item(24, 23).Value = 1.62
$ws.Cells.Item(24, 24).Value = 2.2
$ws.Cells.Item(24, 25).Value = 2.2
$ws.Cells.Item(24, 26).Value = 1.62
$ws.Cells.Item(24, 27).Value = 6.5
$ws.Cells.Item(24, 28).Value = 12
$ws.Cells.Item(24, 29).Value = 12
$ws.Cells.Item(24, 33).Value = 6
$ws.Cells.Item(24, 35).Value = 19
$ws.Cells.Item(24, 40).Value = 11
$ws.Cells.Item(24, 41).Value = 29
$ws.Cells.Item(24, 44).Value = 2
$ws.Cells.Item(24, 45).Value = 1.85
$ws.Cells.Item(25, 7).Value = 3.6
$ws.Cells.Item(25, 8).Value = 3.3
$ws.Cells.Item(25, 9).Value = 2.1
$ws.Cells.Item(25, 10).Value = 4
$ws.Cells.Item(25, 28).Value = 17
$ws.Cells.Item(25, 31).Value = 29
$ws.Cells.Item(25, 39).Value = 9.5
$ws.Cells.Item(25, 41).Value = 19
$ws.Cells.Item(26, 8).Value = 5
$ws.Cells.Item(26, 11).Value = 2.5
$ws.Cells.Item(26, 15).Value = 1.18
$ws.Cells.Item(26, 16).Value = 4.5
$ws.Cells.Item(26, 17).Value = 1.62
$ws.Cells.Item(26, 18).Value = 2.25
$ws.Cells.Item(26, 21).Value = 2.5
$ws.Cells.Item(26, 22).Value = 1.5
$ws.Cells.Item(26, 23).Value = 1.3
$ws.Cells.Item(26, 24).Value = 3.4
$ws.Cells.Item(26, 25).Value = 1.95
$ws.Cells.Item(26, 26).Value = 1.8
$ws.Cells.Item(26, 29).Value = 21
$ws.Cells.Item(26, 30).Value = 81
$ws.Cells.Item(26, 33).Value = 13
$ws.Cells.Item(26, 36).Value = 51
$ws.Cells.Item(26, 37).Value = 301
$ws.Cells.Item(26, 38).Value = 7.5
$ws.Cells.Item(26, 39).Value = 7
$ws.Cells.Item(26, 41).Value = 9
$ws.Cells.Item(26, 42).Value = 11
$ws.Cells.Item(26, 43).Value = 26
$ws.Cells.Item(27, 16).Value = 3.85
$ws.Cells.Item(28, 37).Value = 1250
$ws.Cells.Item(29, 7).Value = 6.5
$ws.Cells.Item(29, 8).Value = 4.33
$ws.Cells.Item(29, 9).Value = 1.42
$ws.Cells.Item(29, 10).Value = 6
$ws.Cells.Item(29, 11).Value = 2.75
$ws.Cells.Item(29, 12).Value = 1.83
$ws.Cells.Item(29, 25).Value = 1.57
$ws.Cells.Item(29, 26).Value = 2.25
$ws.Cells.Item(29, 27).Value = 26
$ws.Cells.Item(29, 29).Value = 21
$ws.Cells.Item(29, 30).Value = 81
$ws.Cells.Item(29, 32).Value = 41
$ws.Cells.Item(29, 34).Value = 9.5
$ws.Cells.Item(29, 36).Value = 41
$ws.Cells.Item(29, 37).Value = 126
$ws.Cells.Item(29, 38).Value = 11
$ws.Cells.Item(29, 39).Value = 9.5
$ws.Cells.Item(29, 41).Value = 11
